$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "MEC-2B-Des. Maq. Cad."
$ws.Range("C2").Value = "MEC-2B-Des. Maq. Cad."
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = "-"
$ws.Range("F2").Value = "MCT-3A-CAM"

$ws.Range("B3").Value = "MCT-3A-CAM"
$ws.Range("E3").Value = "-"
$ws.Range("F3").Value = "MEC-3B-Metrologia 2"

$ws.Range("B4").Value = "MCT-3A-CAM"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "-"
$ws.Range("F4").Value = "MEC-3B-Metrologia 2"

$ws.Range("B6").Value = "MEC-2B-Des. Maq. Cad."
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "-"
$ws.Range("F6").Value = "MCT-3A-CAM"

$ws.Range("B7").Value = "MEC-2B-Des. Maq. Cad."
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = "-"
$ws.Range("F7").Value = "MEC-3B-Metrologia 2"

$ws.Range("B8").Value = "MEC-2B-Des. Maq. Cad."
$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = "MEC-2B-Des. Maq. Cad."
$ws.Range("F8").Value = "MEC-3B-Metrologia 2"

$ws.Range("F10").Value = "MEC-2A-Des. Maq. Cad."

$ws.Range("B11").Value = "-"
$ws.Range("F11").Value = "MEC-2A-Des. Maq. Cad."

$ws.Range("B12").Value = "MEC-2A-Ajustagem"
$ws.Range("F12").Value = "MEC-2A-Des. Maq. Cad."

$ws.Range("B14").Value = "MEC-2A-Ajustagem"
$ws.Range("F14").Value = "MEC-2A-Des. Maq. Cad."

$ws.Range("B15").Value = "MEC-2A-Ajustagem"
$ws.Range("F15").Value = "MEC-2A-Des. Maq. Cad."

$ws.Range("B16").Value = "MEC-2A-Ajustagem"
$ws.Range("F16").Value = "MEC-2A-Des. Maq. Cad."
